$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (97) of data below the existing table, continuing the
# monthly date series. Column A keeps the same date formatting as the row
# above it (copy the format from A96 so no new style is created), then the
# date/Services/Goods values are written into A97:C97.
$ws.Range("A96").Copy()
$ws.Range("A97").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A97").Value = 45627
$ws.Range("B97").Value = 0.223799942512216
$ws.Range("C97").Value = 0.0775162849338096

$excel.CutCopyMode = 0
